$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Value corrections to existing cells (rows 132,135,136) ---
$ws.Range("AM132").Value = 13.3799812
$ws.Range("H135").Value = 12.7323076
$ws.Range("C136").Value = 19.8155376
$ws.Range("G136").Value = 11.4576207
$ws.Range("H136").Value = 12.5175809
$ws.Range("M136").Value = 15.2638941
$ws.Range("U136").Value = 14.5807514
$ws.Range("AP136").Value = 12.0252072

# --- New AR cells added to rows 131-133 ---
$ws.Range("AR131").Value = 13.3621726
$ws.Range("AR132").Value = 14.6226415
$ws.Range("AR133").Value = 11.6381778

# --- Fill in full data rows for 137 and 138 (previously date-only) ---
# Row 137
$ws.Range("B137").Value = 15.6149388
$ws.Range("C137").Value = 19.5283203
$ws.Range("D137").Value = 18.442629
$ws.Range("F137").Value = 18.8023214
$ws.Range("G137").Value = 11.4473773
$ws.Range("H137").Value = 12.4193695
$ws.Range("I137").Value = 11.8438359
$ws.Range("J137").Value = 14.2045455
$ws.Range("K137").Value = 12.0906236
$ws.Range("L137").Value = 12.2196539
$ws.Range("M137").Value = 15.4590936
$ws.Range("O137").Value = 7.6448498
$ws.Range("P137").Value = 18.8169343
$ws.Range("Q137").Value = 14.2532945
$ws.Range("R137").Value = 13.5544847
$ws.Range("S137").Value = 18.2489189
$ws.Range("T137").Value = 14.0046708
$ws.Range("U137").Value = 14.3636571
$ws.Range("V137").Value = 18.6298761
$ws.Range("W137").Value = 13.1804951
$ws.Range("X137").Value = 14.0895147
$ws.Range("Y137").Value = 8.690828399999999
$ws.Range("Z137").Value = 11.2017161
$ws.Range("AA137").Value = 14.6757238
$ws.Range("AB137").Value = 14.3026208
$ws.Range("AD137").Value = 20.7763476
$ws.Range("AE137").Value = 10.2391625
$ws.Range("AF137").Value = 14.8888856
$ws.Range("AG137").Value = 15.7382285
$ws.Range("AH137").Value = 20.285459
$ws.Range("AI137").Value = 11.1610793
$ws.Range("AJ137").Value = 13.4696817
$ws.Range("AK137").Value = 14.1425695
$ws.Range("AL137").Value = 12.9820137
$ws.Range("AM137").Value = 12.4366508
$ws.Range("AN137").Value = 13.2853985
$ws.Range("AO137").Value = 15.7899192
$ws.Range("AP137").Value = 12.0123571
$ws.Range("AQ137").Value = 11.3780708
$ws.Range("AS137").Value = 11.8574575
$ws.Range("AT137").Value = 17.939465
$ws.Range("AU137").Value = 20.6017187
$ws.Range("AV137").Value = 13.75597
$ws.Range("AW137").Value = 16.7361612
$ws.Range("AX137").Value = 18.5291131
$ws.Range("AY137").Value = 13.2754233
$ws.Range("BA137").Value = 7.4195435
$ws.Range("BB137").Value = 13.4423746
$ws.Range("BC137").Value = 13.3799606
$ws.Range("BD137").Value = 13.2352726
$ws.Range("BE137").Value = 13.8654836

# Row 138
$ws.Range("B138").Value = 15.0680493
$ws.Range("C138").Value = 19.6991256
$ws.Range("D138").Value = 18.7396905
$ws.Range("F138").Value = 19.0074847
$ws.Range("G138").Value = 11.4626189
$ws.Range("H138").Value = 12.5664047
$ws.Range("I138").Value = 12.0011983
$ws.Range("J138").Value = 12.8846154
$ws.Range("K138").Value = 12.2835005
$ws.Range("L138").Value = 12.3919952
$ws.Range("M138").Value = 15.8000529
$ws.Range("O138").Value = 8.358433700000001
$ws.Range("P138").Value = 19.0493224
$ws.Range("Q138").Value = 13.9140686
$ws.Range("R138").Value = 13.4160883
$ws.Range("S138").Value = 18.4431473
$ws.Range("T138").Value = 13.8102217
$ws.Range("U138").Value = 14.8422778
$ws.Range("V138").Value = 18.9529513
$ws.Range("W138").Value = 13.2978148
$ws.Range("X138").Value = 13.8270077
$ws.Range("Y138").Value = 9.026920199999999
$ws.Range("Z138").Value = 11.0033129
$ws.Range("AA138").Value = 14.7446173
$ws.Range("AB138").Value = 14.3728993
$ws.Range("AD138").Value = 21.8431047
$ws.Range("AE138").Value = 10.3112981
$ws.Range("AF138").Value = 15.225745
$ws.Range("AG138").Value = 16.4131367
$ws.Range("AH138").Value = 20.4613024
$ws.Range("AI138").Value = 11.4009534
$ws.Range("AJ138").Value = 13.1389045
$ws.Range("AK138").Value = 14.1548762
$ws.Range("AL138").Value = 12.9900376
$ws.Range("AM138").Value = 12.5457224
$ws.Range("AN138").Value = 13.3639204
$ws.Range("AO138").Value = 16.1910891
$ws.Range("AP138").Value = 12.0042083
$ws.Range("AQ138").Value = 11.3281086
$ws.Range("AS138").Value = 11.6375364
$ws.Range("AT138").Value = 18.4908205
$ws.Range("AU138").Value = 20.4183064
$ws.Range("AV138").Value = 13.9301211
$ws.Range("AW138").Value = 17.2901219
$ws.Range("AX138").Value = 18.8928171
$ws.Range("AY138").Value = 13.1343996
$ws.Range("BA138").Value = 8.112098100000001
$ws.Range("BB138").Value = 13.5410959
$ws.Range("BC138").Value = 13.5468238
$ws.Range("BD138").Value = 12.4645669
$ws.Range("BE138").Value = 14.346529

# --- New date labels appended to the date column (rows 139-144) ---
$ws.Range("A139").Value = "17 06 2020"
$ws.Range("A140").Value = "18 06 2020"
$ws.Range("A141").Value = "19 06 2020"
$ws.Range("A142").Value = "20 06 2020"
$ws.Range("A143").Value = "21 06 2020"
$ws.Range("A144").Value = "22 06 2020"

# --- New data rows 139-142 ---
# Row 139
$ws.Range("B139").Value = 15.6052963
$ws.Range("C139").Value = 20.3457794
$ws.Range("D139").Value = 18.9336715
$ws.Range("F139").Value = 19.5502075
$ws.Range("G139").Value = 11.6152557
$ws.Range("H139").Value = 12.5942755
$ws.Range("I139").Value = 12.0476974
$ws.Range("J139").Value = 13.6822194
$ws.Range("K139").Value = 12.1043246
$ws.Range("L139").Value = 12.9842945
$ws.Range("M139").Value = 15.9223064
$ws.Range("O139").Value = 8.7925852
$ws.Range("P139").Value = 19.1730414
$ws.Range("Q139").Value = 13.6903944
$ws.Range("R139").Value = 13.4221344
$ws.Range("S139").Value = 18.7086104
$ws.Range("T139").Value = 14.41682
$ws.Range("U139").Value = 15.1817082
$ws.Range("V139").Value = 19.1760074
$ws.Range("W139").Value = 12.9771405
$ws.Range("X139").Value = 13.5839357
$ws.Range("Y139").Value = 8.8306354
$ws.Range("Z139").Value = 10.8661799
$ws.Range("AA139").Value = 14.9491736
$ws.Range("AB139").Value = 14.6663431
$ws.Range("AD139").Value = 22.2423731
$ws.Range("AE139").Value = 10.7033158
$ws.Range("AF139").Value = 15.5663568
$ws.Range("AG139").Value = 16.6079784
$ws.Range("AH139").Value = 19.7499065
$ws.Range("AI139").Value = 11.2616311
$ws.Range("AJ139").Value = 13.1114203
$ws.Range("AK139").Value = 14.437439
$ws.Range("AL139").Value = 13.1227802
$ws.Range("AM139").Value = 12.2138494
$ws.Range("AN139").Value = 13.4502023
$ws.Range("AO139").Value = 16.5600298
$ws.Range("AP139").Value = 11.5623539
$ws.Range("AQ139").Value = 11.2999319
$ws.Range("AS139").Value = 11.6725328
$ws.Range("AT139").Value = 18.6366866
$ws.Range("AU139").Value = 19.9467344
$ws.Range("AV139").Value = 14.1671142
$ws.Range("AW139").Value = 18.3555256
$ws.Range("AX139").Value = 19.0871463
$ws.Range("AY139").Value = 13.2600756
$ws.Range("BA139").Value = 8.133823700000001
$ws.Range("BB139").Value = 13.400099
$ws.Range("BC139").Value = 13.918514
$ws.Range("BD139").Value = 12.9388228
$ws.Range("BE139").Value = 14.6572575

# Row 140
$ws.Range("B140").Value = 15.7637076
$ws.Range("C140").Value = 20.3321878
$ws.Range("D140").Value = 19.1347461
$ws.Range("F140").Value = 19.8888276
$ws.Range("G140").Value = 12.0209684
$ws.Range("H140").Value = 12.8524218
$ws.Range("I140").Value = 12.0209689
$ws.Range("J140").Value = 13.4013605
$ws.Range("K140").Value = 12.3030584
$ws.Range("L140").Value = 13.4168713
$ws.Range("M140").Value = 16.1857632
$ws.Range("O140").Value = 8.2910751
$ws.Range("P140").Value = 18.3994708
$ws.Range("Q140").Value = 14.0356044
$ws.Range("R140").Value = 13.499991
$ws.Range("S140").Value = 18.4196807
$ws.Range("T140").Value = 14.6683156
$ws.Range("U140").Value = 15.0604059
$ws.Range("V140").Value = 19.4530328
$ws.Range("W140").Value = 12.6542629
$ws.Range("X140").Value = 13.6196414
$ws.Range("Y140").Value = 9.153005500000001
$ws.Range("Z140").Value = 10.8909042
$ws.Range("AA140").Value = 15.0131008
$ws.Range("AB140").Value = 14.8439265
$ws.Range("AD140").Value = 22.1329414
$ws.Range("AE140").Value = 10.691112
$ws.Range("AF140").Value = 15.5337024
$ws.Range("AG140").Value = 15.9935856
$ws.Range("AH140").Value = 19.5951362
$ws.Range("AI140").Value = 11.1190818
$ws.Range("AJ140").Value = 12.837067
$ws.Range("AK140").Value = 14.9511402
$ws.Range("AL140").Value = 13.6608688
$ws.Range("AM140").Value = 12.2546746
$ws.Range("AN140").Value = 13.0820181
$ws.Range("AO140").Value = 17.6449543
$ws.Range("AP140").Value = 11.7066135
$ws.Range("AQ140").Value = 11.2940544
$ws.Range("AS140").Value = 11.7879966
$ws.Range("AT140").Value = 19.3322409
$ws.Range("AU140").Value = 21.1851886
$ws.Range("AV140").Value = 14.5164984
$ws.Range("AW140").Value = 19.2173408
$ws.Range("AX140").Value = 18.9721588
$ws.Range("AY140").Value = 13.6086623
$ws.Range("BA140").Value = 9.355458199999999
$ws.Range("BB140").Value = 13.3346916
$ws.Range("BC140").Value = 13.5697173
$ws.Range("BD140").Value = 13.2590809
$ws.Range("BE140").Value = 14.8673203

# Row 141
$ws.Range("B141").Value = 15.6911142
$ws.Range("C141").Value = 20.5188878
$ws.Range("D141").Value = 19.4440121
$ws.Range("F141").Value = 20.4707148
$ws.Range("G141").Value = 12.0936925
$ws.Range("H141").Value = 13.0305609
$ws.Range("I141").Value = 12.1035387
$ws.Range("J141").Value = 14.4787645
$ws.Range("K141").Value = 12.8211462
$ws.Range("L141").Value = 13.8231764
$ws.Range("M141").Value = 16.5580787
$ws.Range("O141").Value = 9.0078329
$ws.Range("P141").Value = 19.2694904
$ws.Range("Q141").Value = 14.5680868
$ws.Range("R141").Value = 13.3354445
$ws.Range("S141").Value = 18.6318073
$ws.Range("T141").Value = 15.2372502
$ws.Range("U141").Value = 15.4480315
$ws.Range("V141").Value = 19.9544709
$ws.Range("W141").Value = 12.3700826
$ws.Range("X141").Value = 13.9466109
$ws.Range("Y141").Value = 8.791691
$ws.Range("Z141").Value = 10.6963159
$ws.Range("AA141").Value = 15.139835
$ws.Range("AB141").Value = 15.4080429
$ws.Range("AD141").Value = 22.1583041
$ws.Range("AE141").Value = 11.9398349
$ws.Range("AF141").Value = 15.6684567
$ws.Range("AG141").Value = 16.0921332
$ws.Range("AH141").Value = 18.9478719
$ws.Range("AI141").Value = 11.3308801
$ws.Range("AJ141").Value = 12.8977338
$ws.Range("AK141").Value = 14.2519614
$ws.Range("AL141").Value = 13.6566247
$ws.Range("AM141").Value = 12.1808948
$ws.Range("AN141").Value = 13.0037343
$ws.Range("AO141").Value = 18.0259254
$ws.Range("AP141").Value = 11.7626904
$ws.Range("AQ141").Value = 11.4475768
$ws.Range("AS141").Value = 11.8889121
$ws.Range("AT141").Value = 20.2169335
$ws.Range("AU141").Value = 21.3479031
$ws.Range("AV141").Value = 14.6958387
$ws.Range("AW141").Value = 20.1087938
$ws.Range("AX141").Value = 19.1362719
$ws.Range("AY141").Value = 13.4989403
$ws.Range("BA141").Value = 9.2434162
$ws.Range("BB141").Value = 13.3152779
$ws.Range("BC141").Value = 13.7299775
$ws.Range("BD141").Value = 13.1334905
$ws.Range("BE141").Value = 14.1792196

# Row 142
$ws.Range("B142").Value = 15.6647808
$ws.Range("C142").Value = 20.760556
$ws.Range("D142").Value = 20.2102959
$ws.Range("F142").Value = 21.1602242
$ws.Range("G142").Value = 12.2418795
$ws.Range("H142").Value = 13.2127952
$ws.Range("I142").Value = 11.8537964
$ws.Range("J142").Value = 13.5994587
$ws.Range("K142").Value = 13.1284916
$ws.Range("L142").Value = 14.4759983
$ws.Range("M142").Value = 16.6395388
$ws.Range("O142").Value = 8.985429
$ws.Range("P142").Value = 19.5761948
$ws.Range("Q142").Value = 14.9041295
$ws.Range("R142").Value = 13.3006391
$ws.Range("S142").Value = 18.6298605
$ws.Range("T142").Value = 15.4976721
$ws.Range("U142").Value = 15.2144724
$ws.Range("V142").Value = 19.9447108
$ws.Range("W142").Value = 12.1804351
$ws.Range("X142").Value = 13.6294858
$ws.Range("Y142").Value = 9.172030400000001
$ws.Range("Z142").Value = 10.7757837
$ws.Range("AA142").Value = 14.9163279
$ws.Range("AB142").Value = 15.1062203
$ws.Range("AD142").Value = 22.8177464
$ws.Range("AE142").Value = 11.1284749
$ws.Range("AF142").Value = 16.1250255
$ws.Range("AG142").Value = 17.2779457
$ws.Range("AH142").Value = 18.8760784
$ws.Range("AI142").Value = 11.112898
$ws.Range("AJ142").Value = 12.6808796
$ws.Range("AK142").Value = 13.8305469
$ws.Range("AL142").Value = 13.7830959
$ws.Range("AM142").Value = 12.1841492
$ws.Range("AN142").Value = 12.7262508
$ws.Range("AO142").Value = 19.1477569
$ws.Range("AP142").Value = 11.8581664
$ws.Range("AQ142").Value = 11.3918073
$ws.Range("AS142").Value = 11.8782537
$ws.Range("AT142").Value = 20.9100572
$ws.Range("AU142").Value = 20.0406169
$ws.Range("AV142").Value = 15.0080442
$ws.Range("AW142").Value = 21.0346173
$ws.Range("AX142").Value = 18.9995117
$ws.Range("AY142").Value = 13.6588603
$ws.Range("BA142").Value = 9.3976778
$ws.Range("BB142").Value = 13.2376404
$ws.Range("BC142").Value = 13.7470551
$ws.Range("BD142").Value = 13.0761004
$ws.Range("BE142").Value = 14.7427201
